$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 200
$ws.Range("B3").Value = 720
$ws.Range("B5").Value = 34
$ws.Range("B6").Value = 183
$ws.Range("B8").Value = 300
